$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: header row (First Name, Last name, age, gender, usename, phone, skype)
$ws.Range("B38").Value = "First Name"
$ws.Range("C38").Value = "Last name"
$ws.Range("D38").Value = "age"
$ws.Range("E38").Value = "gender"
$ws.Range("F38").Value = "usename"
$ws.Range("G38").Value = "phone"
$ws.Range("H38").Value = "skype"

# Row 39: data row
$ws.Range("B39").Value = "Dmytro"
$ws.Range("C39").Value = "Mula"
$ws.Range("D39").Value = 22
$ws.Range("E39").Value = "male"
$ws.Range("F39").Value = "skip"
$ws.Range("G39").Value = 89123123
$ws.Range("H39").Value = "dimon.mula"

# Row 42: header row (same as row 38)
$ws.Range("B42").Value = "First Name"
$ws.Range("C42").Value = "Last name"
$ws.Range("D42").Value = "age"
$ws.Range("E42").Value = "gender"
$ws.Range("F42").Value = "usename"
$ws.Range("G42").Value = "phone"
$ws.Range("H42").Value = "skype"

# Row 47: data row (same as row 39)
$ws.Range("B47").Value = "Dmytro"
$ws.Range("C47").Value = "Mula"
$ws.Range("D47").Value = 22
$ws.Range("E47").Value = "male"
$ws.Range("F47").Value = "skip"
$ws.Range("G47").Value = 89123123
$ws.Range("H47").Value = "dimon.mula"

# Rows 50-56: header/value pairs laid out as columns, B = label, F = value
$ws.Range("B50").Value = "First Name"
$ws.Range("F50").Value = "Dmytro"

$ws.Range("B51").Value = "Last name"
$ws.Range("F51").Value = "Mula"

$ws.Range("B52").Value = "age"
$ws.Range("F52").Value = 22

$ws.Range("B53").Value = "gender"
$ws.Range("F53").Value = "male"

$ws.Range("B54").Value = "usename"
$ws.Range("F54").Value = "skip"

$ws.Range("B55").Value = "phone"
$ws.Range("F55").Value = 89123123

$ws.Range("B56").Value = "skype"
$ws.Range("F56").Value = "dimon.mula"

# Column G now holds data (phone numbers), give it a custom width like column B
$ws.Columns.Item(7).ColumnWidth = 14.5

# Update the view to reflect the new selection/scroll position seen in the target workbook
$excel.Goto($ws.Range("A32"), $true)
$ws.Range("G55").Select()
